$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "thisletter" values for two rows
$ws.Range("B29").Value = "J"
$ws.Range("B59").Value = "Y"

# Update the current selection
$ws.Range("G3").Select()
